$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 47; this shifts the existing rows 47:71 down to 48:72
# (matching the diff's net effect of a new weekly record being prepended to
# this block and every later record shifting down by one row).
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row 47 with the new weekly record.
$ws.Range("A47").Value = 10
$ws.Range("B47").Value = "Vega Modelo de Temuco"
$ws.Range("C47").Value = "La Araucanía"
$ws.Range("D47").Value = 44767
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 300000001
$ws.Range("G47").Value = "Rabanito"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 50
$ws.Range("K47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("M47").Value = 8000
$ws.Range("N47").Value = "$/docena de paquetes"
$ws.Range("O47").Value = "Provincia de Cautín"
$ws.Range("P47").Value = 667
$ws.Range("Q47").Value = 12
$ws.Range("R47").Value = "Hortaliza"
